$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3:C5").NumberFormat = "@"

# Row 3: KHALID VAVA
$ws.Range("A3").Value = "KHALID VAVA "
$ws.Range("B3").Value = "O3546845"
$ws.Range("C3").Value = "321564613641864613156486"
$ws.Range("D3").Value = "AG 5"
$ws.Range("E3").Value = "bmce"
$ws.Range("F3").Value = "Supervision"
$ws.Range("G3").Value = "120/SUP 2"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 18000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 18000

# Row 4: FATIMA FAFA
$ws.Range("A4").Value = "FATIMA FAFA"
$ws.Range("B4").Value = "L3659652"
$ws.Range("C4").Value = "246848931356984893231321"
$ws.Range("D4").Value = "CASA 556"
$ws.Range("E4").Value = "BP"
$ws.Range("F4").Value = "Supervision"
$ws.Range("G4").Value = "120/SUP 2"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 12000
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 11400

# Row 5: IMANE FAFA
$ws.Range("A5").Value = "IMANE FAFA"
$ws.Range("B5").Value = "K6546841"
$ws.Range("C5").Value = "313215156145641564165411"
$ws.Range("D5").Value = "AGG55"
$ws.Range("E5").Value = "BMCE"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "001/CASA NORD"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 48000
$ws.Range("J5").Value = 5400
$ws.Range("K5").Value = 42600
